# Weekly fruit/vegetable price update: insert a new record for the week
# right after the existing row 241, shifting subsequent rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 241; rows 241..260 shift down to 242..261.
$ws.Rows.Item(241).Insert()

# Populate the new row 241 with the new weekly price entry. Most of the
# fields mirror what used to be in (old) row 241 — only the date and the
# price columns change for the new week's observation.
$ws.Range("A241").Value = 8
$ws.Range("B241").Value = "Terminal La Palmera de La Serena"
$ws.Range("C241").Value = "Coquimbo"
$ws.Range("D241").Value = 44783
$ws.Range("E241").Value = 4
$ws.Range("F241").Value = 100112031
$ws.Range("G241").Value = "Poroto verde"
$ws.Range("H241").Value = "Magnum"
$ws.Range("I241").Value = "Primera"
$ws.Range("J241").Value = 500
$ws.Range("K241").Value = 34500
$ws.Range("L241").Value = 35000
$ws.Range("M241").Value = 34750
$ws.Range("N241").Value = "$/malla 25 kilos"
$ws.Range("O241").Value = "Perú"
$ws.Range("P241").Value = 1390
$ws.Range("Q241").Value = 25
$ws.Range("R241").Value = "Hortaliza"
